$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") values for rows 2-8 move from 45207 to 45208 (i.e. +1 day)
foreach ($row in 2..8) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45208
}
